# Auto-generated edit script: update crypto price/volume table
# to reflect the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.404.57'
$ws.Range("E2").Value = '  +7.12%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.545.44'
$ws.Range("E3").Value = '  +10.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '191.25'
$ws.Range("E5").Value = '  +9.82%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '554.74'
$ws.Range("E6").Value = '  +4.36%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.540.22'
$ws.Range("E7").Value = '  +9.96%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.608'
$ws.Range("E8").Value = '  +2.14%  '

# Row 9
$ws.Range("E9").Value = '  -0.06%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.636'
$ws.Range("E10").Value = '  +4.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("E11").Value = '  +17.13%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.30'
$ws.Range("E12").Value = '  +2.93%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  +7.70%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.43'
$ws.Range("E14").Value = '  +3.44%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.079.20'
$ws.Range("E15").Value = '  +9.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.520.96'
$ws.Range("E16").Value = '  +9.47%  '

# Row 17
$ws.Range("E17").Value = '  +3.74%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.241.27'
$ws.Range("E18").Value = '  +7.20%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.28'
$ws.Range("E19").Value = '  +5.17%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.99'
$ws.Range("E20").Value = '  +8.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.996'
$ws.Range("E21").Value = '  +2.63%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '433.64'
$ws.Range("E22").Value = '  +18.09%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '85.52'
$ws.Range("E23").Value = '  +5.24%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.92'
$ws.Range("E24").Value = '  +3.74%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.16'
$ws.Range("E25").Value = '  +6.34%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.17'
$ws.Range("E26").Value = '  -0.64%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.92'
$ws.Range("E27").Value = '  +9.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.08'
$ws.Range("E28").Value = '  +6.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.03'
$ws.Range("E29").Value = '  +10.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.43'
$ws.Range("E30").Value = '  +6.45%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '644.19'
$ws.Range("E31").Value = '  -0.20%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.71'
$ws.Range("E32").Value = '  +2.66%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.79'
$ws.Range("E33").Value = '  +3.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.111'
$ws.Range("E34").Value = '  +4.77%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.71'
$ws.Range("E35").Value = '  +5.15%  '

# Row 36
$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0835'
$ws.Range("E36").Value = '  +16.20%  '

# Row 37
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.65'
$ws.Range("E37").Value = '  +4.42%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.12%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.144'
$ws.Range("E39").Value = '  +16.13%  '

# Row 40
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.393'
$ws.Range("E40").Value = '  +3.76%  '

# Row 41
$ws.Range("E41").Value = '  +13.15%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.995'
$ws.Range("E42").Value = '  -0.21%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.037.81'
$ws.Range("E43").Value = '  +5.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.65'
$ws.Range("E44").Value = '  +5.06%  '

# Row 45
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.88'
$ws.Range("E45").Value = '  +10.66%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.87'
$ws.Range("E46").Value = '  +6.77%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  +9.94%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0420'
$ws.Range("E48").Value = '  +6.40%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.131'
$ws.Range("E49").Value = '  +5.20%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.73'
$ws.Range("E50").Value = '  +12.72%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '140.42'
$ws.Range("E51").Value = '  +4.22%  '
